# Update values in column A (and a few in column E) for the "AE" combination
# RandomForest imputed result data. This mirrors an updated algorithm run
# where some imputed values differ slightly from the previous ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4"   = -20.95219999999999
    "A7"   = -20.46869999999997
    "E7"   = 15.3553
    "E15"  = 16.314
    "A16"  = -21.71689999999999
    "E21"  = 17.1187
    "E22"  = 17.027
    "E23"  = 16.18009999999998
    "A28"  = -21.93150000000001
    "A29"  = -21.27519999999998
    "A32"  = -21.1695
    "E34"  = 17.3585
    "A40"  = -20.16569999999999
    "E43"  = 17.49470000000001
    "E45"  = 16.3954
    "E50"  = 16.2108
    "E51"  = 17.28500000000001
    "A52"  = -22.1758
    "A57"  = -22.28290000000001
    "A66"  = -21.4987
    "E66"  = 17.12910000000002
    "E67"  = 17.09330000000001
    "E79"  = 18.45300000000002
    "E84"  = 16.5426
    "E92"  = 18.44610000000002
    "E97"  = 16.66289999999999
    "A100" = -21.9452
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
